$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits near the
#    top of the document (right after the phone/email line's trailing
#    tab run). It will be re-created later at the new edit location.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the text we need to touch: "...jQuery, beginner React.js"
#    becomes "...jQuery, beginner at React.js".
# ------------------------------------------------------------------
$reactRng = $d.Content
$null = $reactRng.Find.Execute("React", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$reactStart = $reactRng.Start
$reactEnd = $reactRng.End

$jqueryRng = $d.Content
$null = $jqueryRng.Find.Execute("jQuery", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jqueryEnd = $jqueryRng.End

# ------------------------------------------------------------------
# 3. Pin down every run boundary we want to keep with throw-away
#    bookmarks *before* inserting text -- Word (and this host) merges
#    same-formatted runs across an edited paragraph unless a bookmark
#    (or similar hard stop) protects the seam.
# ------------------------------------------------------------------
$mJQueryEnd = $d.Range($jqueryEnd, $jqueryEnd)
$d.Bookmarks.Add("ZZSPLIT_JQUERY", $mJQueryEnd) | Out-Null

$mReactStart = $d.Range($reactStart, $reactStart)
$d.Bookmarks.Add("ZZSPLIT_REACT", $mReactStart) | Out-Null

$mReactEnd = $d.Range($reactEnd, $reactEnd)
$d.Bookmarks.Add("ZZSPLIT_JS", $mReactEnd) | Out-Null

# ------------------------------------------------------------------
# 4. Type the new word "at " immediately before "React".
# ------------------------------------------------------------------
$insertionPoint = $d.Range($reactStart, $reactStart)
$insertionPoint.InsertAfter("at ")

# ------------------------------------------------------------------
# 5. Split "at " off from ", beginner " into its own run by dropping
#    another boundary bookmark exactly between them.
# ------------------------------------------------------------------
$reactBm = $d.Bookmarks("ZZSPLIT_REACT")
$atStart = $reactBm.Start - 3
$mAtStart = $d.Range($atStart, $atStart)
$d.Bookmarks.Add("ZZSPLIT_AT", $mAtStart) | Out-Null

# ------------------------------------------------------------------
# 6. Drop the helper bookmarks we no longer need; the run splits they
#    created stay intact once the paragraph has been normalized.
# ------------------------------------------------------------------
$d.Bookmarks("ZZSPLIT_JQUERY").Delete()
$d.Bookmarks("ZZSPLIT_JS").Delete()

# ------------------------------------------------------------------
# 7. Re-create "_GoBack" at the new edit location -- right between
#    "at " and "React" -- matching where Word leaves it after the
#    last piece of typed text, and drop the temporary marker there.
# ------------------------------------------------------------------
$finalBm = $d.Bookmarks("ZZSPLIT_REACT")
$goBackRange = $d.Range($finalBm.Start, $finalBm.Start)
$finalBm.Delete()
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$d.Bookmarks("ZZSPLIT_AT").Delete()

Write-Output $d.Paragraphs.Item(21).Range.Text
